# Weekly update: insert a new most-recent-week price record for
# "Vega Modelo de Temuco" / Pomelo, pushing the existing history rows
# (397-413) down by one (to 398-414).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 397; everything that was
# in rows 397:413 shifts down to 398:414 (dimension grows to T414).
$ws.Rows.Item(397).Insert()

# Populate the newly inserted row 397 with this week's data.
$ws.Range("A397").Value = 10
$ws.Range("B397").Value = "Vega Modelo de Temuco"
$ws.Range("C397").Value = "La Araucanía"
$ws.Range("D397").Value = 45041
$ws.Range("E397").Value = 9
$ws.Range("F397").Value = "Fruta"
$ws.Range("G397").Value = 100102
$ws.Range("H397").Value = "Cítricos"
$ws.Range("I397").Value = 100102006
$ws.Range("J397").Value = "Pomelo"
$ws.Range("K397").Value = "Start Ruby"
$ws.Range("L397").Value = "Primera"
$ws.Range("M397").Value = 65
$ws.Range("N397").Value = 15000
$ws.Range("O397").Value = 15000
$ws.Range("P397").Value = 15000
$ws.Range("Q397").Value = "$/bandeja 15 kilos granel"
$ws.Range("R397").Value = "Región de O'Higgins"
$ws.Range("S397").Value = 1000
$ws.Range("T397").Value = 15
